# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

$ws.Range("D3").Value = 92.40000000000001
$ws.Range("D4").Value = 98.59999999999999
$ws.Range("C7").Value = 158
$ws.Range("C8").Value = 272
